# Helplines slide: resize "GoBusiness Licensing" rectangle/textbox,
# update the licensing-helpdesk email to the .com.sg domain, split the
# GoBusiness licence URL into the new gobusiness.gov.sg link, and tidy
# up a few runs that had been split needlessly.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

$EMU = 12700.0
# Shape.Left/Top/Width/Height round-trip through a single-precision
# (points) value, which can truncate by 1 EMU on the way back out.
# Nudge by a hair so the EMU we ask for is the EMU we get.
$EPS = 0.00003

# ---------------------------------------------------------------
# 1. Resize "Rectangle 10" (the shaded background behind the
#    GoBusiness Licensing column) - shift left edge and widen it.
# ---------------------------------------------------------------
$rect10 = $s.Shapes.Item(5)
$rect10.Left = (4178216 / $EMU) + $EPS
$rect10.Width = (3732516 / $EMU) + $EPS

# ---------------------------------------------------------------
# 2. Resize "TextBox 4" (the GoBusiness Licensing contact card) to
#    match, and fix up its text.
# ---------------------------------------------------------------
$tb4 = $s.Shapes.Item(7)
$tb4.Left = (4273911 / $EMU) + $EPS
$tb4.Width = (3543645 / $EMU) + $EPS

$tr4 = $tb4.TextFrame.TextRange

# Email address now ends in .com.sg
$emailPara = $tr4.Paragraphs(6)
$emailPrefix = "Email: "
$emailRun = $tr4.Characters($emailPara.Start + $emailPrefix.Length, $emailPara.Length - $emailPrefix.Length)
$emailRun.Text = "licences-helpdesk@crimsonlogic.com.sg"

# Website now points at the new GoBusiness licence portal. Keep
# "https" as its own run and put the rest ("://...") in a second run,
# mirroring the source edit.
$webPara = $tr4.Paragraphs(7)
$webPrefix = "Website: "
$httpsRun = $tr4.Characters($webPara.Start + $webPrefix.Length, 5)
$httpsRun.Text = "https"
$webPara = $tr4.Paragraphs(7)
$restStart = $webPara.Start + $webPrefix.Length + 5
$restRun = $tr4.Characters($restStart, $webPara.Start + $webPara.Length - $restStart)
$restRun.Text = "://www.gobusiness.gov.sg/licences"

# ---------------------------------------------------------------
# 3. "TextBox 12" (For enquiries relating to GoBusiness Licensing
#    matters...) had a few runs split apart for no reason; merge
#    them back into single runs.
# ---------------------------------------------------------------
$tb12 = $s.Shapes.Item(8)
$tr12 = $tb12.TextFrame.TextRange

# Para 1: "...GoBusiness Licensing" + " " + "matters, e.g.:" -> "...GoBusiness Licensing" + " matters, e.g.:"
$para1 = $tr12.Paragraphs(1)
$prefix1 = "For enquiries relating to GoBusiness Licensing"
$rest1 = $tr12.Characters($para1.Start + $prefix1.Length, $para1.Length - $prefix1.Length)
$rest1.Text = " matters, e.g.:"

# Para 2: "GoBusiness" + " Licensing " + "Login ID" -> "GoBusiness" + " Licensing Login ID"
$para2 = $tr12.Paragraphs(2)
$prefix2 = "GoBusiness"
$rest2 = $tr12.Characters($para2.Start + $prefix2.Length, $para2.Length - $prefix2.Length)
$rest2.Text = " Licensing Login ID"

# Para 3: "forgot/reset GoBusiness" + " Licensing " + "password" -> "forgot/reset GoBusiness" + " Licensing password"
$para3 = $tr12.Paragraphs(3)
$prefix3 = "forgot/reset GoBusiness"
$rest3 = $tr12.Characters($para3.Start + $prefix3.Length, $para3.Length - $prefix3.Length)
$rest3.Text = " Licensing password"
